$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.828.67"
$ws.Range("E2").Value = "  +1.48%  "

$ws.Range("D3").Value = "3.470.31"
$ws.Range("E3").Value = "  +1.74%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "415.01"
$ws.Range("E5").Value = "  +1.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.25"
$ws.Range("E6").Value = "  +1.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  -0.33%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.727"
$ws.Range("E9").Value = "  -0.54%  "

$ws.Range("E10").Value = "  +10.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.56"
$ws.Range("E11").Value = "  -0.26%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000233"
$ws.Range("E12").Value = "  +7.79%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.68"
$ws.Range("E13").Value = "  +5.89%  "

$ws.Range("D14").Value = "4.022.25"
$ws.Range("E14").Value = "  +1.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.57"
$ws.Range("E16").Value = "  -2.79%  "

$ws.Range("D17").Value = "3.491.01"
$ws.Range("E17").Value = "  +2.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.62"
$ws.Range("E18").Value = "  +1.38%  "

$ws.Range("E19").Value = "  -0.60%  "

$ws.Range("D20").Value = "62.741.16"
$ws.Range("E20").Value = "  +1.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "462.12"
$ws.Range("E21").Value = "  +2.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "90.48"
$ws.Range("E22").Value = "  -0.86%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.28"
$ws.Range("E23").Value = "  +2.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.25"
$ws.Range("E24").Value = "  +1.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.72"
$ws.Range("E25").Value = "  +15.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.32"
$ws.Range("E26").Value = "  +1.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.41"
$ws.Range("E27").Value = "  +1.45%  "

$ws.Range("E28").Value = "  +0.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.56"
$ws.Range("E29").Value = "  -0.60%  "

$ws.Range("E30").Value = "  +0.35%  "

$ws.Range("E31").Value = "  -0.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.168"
$ws.Range("E32").Value = "  -1.83%  "

$ws.Range("E33").Value = "  -1.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.96"
$ws.Range("E34").Value = "  -3.93%  "

$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.38"
$ws.Range("E36").Value = "  +8.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0491"
$ws.Range("E37").Value = "  -1.88%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.09"
$ws.Range("E38").Value = "  +5.05%  "

$ws.Range("E39").Value = "  +0.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "149.29"
$ws.Range("E40").Value = "  +4.90%  "

$ws.Range("E41").Value = "  +6.95%  "

$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.134"
$ws.Range("E42").Value = "  -0.10%  "

$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.321"
$ws.Range("E43").Value = "  +1.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.32"
$ws.Range("E44").Value = "  -1.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.39"
$ws.Range("E45").Value = "  +3.61%  "

$ws.Range("E46").Value = "  +3.49%  "

$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.39"
$ws.Range("E47").Value = "  +11.78%  "

$ws.Range("B48").Value = "PEPE"
$ws.Range("C48").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D48").Value = "0.0₃0563"
$ws.Range("E48").Value = "  +34.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.38"
$ws.Range("E49").Value = "  -1.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.39"
$ws.Range("E50").Value = "  +0.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.141"
$ws.Range("E51").Value = "  -2.23%  "

